# update database and change read_price algorithm
# Shift the yearly columns (E..I) one year forward (drop the oldest year,
# add a brand-new year of data) across every block of the "Overview" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the year-label headers (row 8, 16, 25, 34, 42, 51) so the five
#    visible periods move from 1396-1400 to 1397-1401.
# ---------------------------------------------------------------------------
$yearLabels = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)
$yearCols = @("E", "F", "G", "H", "I")
$headerRows = @(8, 16, 25, 34, 42, 51)

foreach ($r in $headerRows) {
    for ($i = 0; $i -lt 5; $i++) {
        $addr = "$($yearCols[$i])$r"
        $ws.Range($addr).Value2 = $yearLabels[$i]
    }
}

# ---------------------------------------------------------------------------
# 2) Helper: shifts a data row's E:I values one column to the left and
#    writes a brand-new value into column I. A value of "-" (string) means
#    "no data" and is written as text; anything else is written as a number.
# ---------------------------------------------------------------------------
function Set-RowValues {
    param($Row, $Values)

    $cols = @("E", "F", "G", "H", "I")
    for ($i = 0; $i -lt 5; $i++) {
        $addr = "$($cols[$i])$Row"
        $ws.Range($addr).Value2 = $Values[$i]
    }
}

# ---------------------------------------------------------------------------
# 3) Apply the new values (previous F:I shifted into E:H, new figure in I)
#    for every data row of every block.
# ---------------------------------------------------------------------------

# Block: مقدار تولید (production quantity) - rows 10,11,12
Set-RowValues 10 @(43127, 42737, 46502, 51249, 60355)
Set-RowValues 11 @("-", "-", "-", 0, 0)
Set-RowValues 12 @(43127, 42737, 46502, 51249, 60355)

# Block: مقدار فروش (sales quantity) - rows 18,19,20,21
Set-RowValues 18 @(43589, 41927, 49328, 52536, 63663)
Set-RowValues 19 @(-631, -437, -253, "-", "-")
Set-RowValues 20 @(0, 0, 0, -271, -266)
Set-RowValues 21 @(42958, 41490, 49075, 52265, 63397)

# Block: مبلغ فروش (sales amount) - rows 27,28,29,30
Set-RowValues 27 @(2732200, 4072381, 8504621, 13752487, 24728679)
Set-RowValues 28 @(-38923, -41115, -36085, "-", "-")
Set-RowValues 29 @(0, 0, 0, -61206, -82738)
Set-RowValues 30 @(2693277, 4031266, 8468536, 13691281, 24645941)

# Block: نرخ فروش (sales rate) - rows 36,37,38
Set-RowValues 36 @(62680952, 97130274, 172409605, 261772632, 388430941)
Set-RowValues 37 @(0, 0, 0, "-", "-")
Set-RowValues 38 @(0, 0, 0, 0, 0)

# Block: مبلغ بهای تمام شده (cost of goods amount) - rows 44,45,46,47
Set-RowValues 44 @(-2728265, -3635503, -7330719, -11372190, -19861474)
Set-RowValues 45 @(39199, 39963, 33487, "-", "-")
Set-RowValues 46 @(0, 0, 0, 52920, 76299)
Set-RowValues 47 @(-2689066, -3595540, -7297232, -11319270, -19785175)

# Block: سود ناخالص (gross profit) - rows 53,54,55,56
Set-RowValues 53 @(3935, 436878, 1173902, 2380297, 4867205)
Set-RowValues 54 @(276, -1152, -2598, "-", "-")
Set-RowValues 55 @("-", 0, 0, -8286, -6439)
Set-RowValues 56 @(4211, 435726, 1171304, 2372011, 4860766)
